$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text: set apostrophe-prefixed value (quote-prefix keeps it text even
# for numeric-looking strings), then reset Style so the cell keeps its
# original (default) style index instead of picking up quotePrefix styling.

$ws.Range("D2").Value = "'27.828.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "'1.870.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'313.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "'0.4826"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("D8").Value = "'0.3812"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("D9").Value = "'0.07372"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").Value = "'0.9407"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "'21.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("D12").Value = "'0.07815"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'1.879.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Value = "'5.492"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "'6.617"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").Value = "'91.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").Value = "'1.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "'0.000008847"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").Value = "'1.012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "'27.836.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "'14.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "'5.120"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "'2.118.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("D24").Value = "'10.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("D25").Value = "'1.949"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "'157.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("D27").Value = "'18.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").Value = "'2.050"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("D29").Value = "'115.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").Value = "'4.979"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "'0.08905"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").Value = "'3.338"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "'1.229"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.93%  "
$ws.Range("D34").Value = "'0.7672"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("D35").Value = "'4.658"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").Value = "'2.737"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").Value = "'1.136"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("D38").Value = "'0.02045"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").Value = "'0.5608"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.26%  "
$ws.Range("D40").Value = "'0.05362"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").Value = "'7.045"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'8.545"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.59%  "
$ws.Range("D44").Value = "'0.1530"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D49").Value = "'1.664"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'68.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("D51").Value = "'0.06117"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "

# Row swaps (B, C, D, E) for rows 45-48 (two coins switched rank order)
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.4874"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.013"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'105.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.96%  "
